$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates (reference period rolled forward to October 2025) ---
$ws.Range("A16").Value = "b. Includes all deaths (both doctor and coroner certified) that occurred and were registered by 31 October 2025."
$ws.Range("A21").Value = "Source: Australian Bureau of Statistics, Deaths due to acute respiratory infections in Australia October 2025"

# --- Row 6 (2024 - RSV): Dec + All revised ---
$ws.Range("M6").Value = 297
$ws.Range("N6").Value = 5105

# --- Row 9 (2025 - COVID-19): Jan, Jun, Jul, Aug, Sep, Oct (newly published) + All revised ---
$ws.Range("B9").Value = 323
$ws.Range("G9").Value = 354
$ws.Range("H9").Value = 359
$ws.Range("I9").Value = 207
$ws.Range("J9").Value = 139
$ws.Range("K9").NumberFormat = "#,##0"
$ws.Range("K9").Value = 38
$ws.Range("N9").Value = 2004

# --- Row 10 (2025 - Influenza): May, Jun, Jul, Aug, Sep, Oct (newly published) + All revised ---
$ws.Range("F10").Value = 79
$ws.Range("G10").Value = 167
$ws.Range("H10").Value = 322
$ws.Range("I10").Value = 286
$ws.Range("J10").Value = 245
$ws.Range("K10").NumberFormat = "#,##0"
$ws.Range("K10").Value = 59
$ws.Range("N10").Value = 1385

# --- Row 11 (2025 - RSV): Jun, Aug, Sep, Oct (newly published) + All revised ---
$ws.Range("G11").Value = 66
$ws.Range("I11").Value = 97
$ws.Range("J11").Value = 77
$ws.Range("K11").NumberFormat = "#,##0"
$ws.Range("K11").Value = 20
$ws.Range("N11").Value = 506
